# updated CB_API and Dash
# Updates a block of computed values (row 2-5) in Sheet1 to reflect the
# refreshed CB_API / Dash pull.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("X2").Value = 0
$ws.Range("Y2").Value = 66.19027777777778
$ws.Range("Z2").Value = 53.53392857142858
$ws.Range("AA2").Value = 55.15363175675675
$ws.Range("AB2").Value = 58.160625

# Row 3
$ws.Range("D3").Value = 49.41402925531915
$ws.Range("E3").Value = 50.22132352941176
$ws.Range("F3").Value = 40.38193359375
$ws.Range("N3").Value = 55.80778301886792
$ws.Range("O3").Value = 56.88100961538461
$ws.Range("P3").Value = 57.00367647058824
$ws.Range("Q3").Value = 79.53749999999999
$ws.Range("R3").Value = 97.04122340425532
$ws.Range("S3").Value = 93.12499999999999
$ws.Range("Y3").Value = 75.4140625
$ws.Range("Z3").Value = 75.78323863636365
$ws.Range("AA3").Value = 46.03629807692307
$ws.Range("AB3").Value = 31.83579545454545
$ws.Range("AD3").Value = 32.32058823529412
$ws.Range("AE3").Value = 45.325
$ws.Range("AF3").Value = 51.93913043478261
$ws.Range("AG3").Value = 43.57790948275862
$ws.Range("AH3").Value = 40.37410714285714

# Row 4
$ws.Range("D4").Value = 44.16201923076923
$ws.Range("F4").Value = 41.67732558139534
$ws.Range("H4").Value = 99.328125
$ws.Range("X4").Value = 39.261328125
$ws.Range("Y4").Value = 38.446875

# Row 5
$ws.Range("T5").Value = 72.01875
$ws.Range("U5").Value = 68.60937499999999
$ws.Range("V5").Value = 38.984375
